# Update "想去人数" (want-to-go count) figures in column F for a handful of
# rows on the "展览" and "全部类型" worksheets. Each value is bumped up by
# the amounts shown below (most +1, one +3).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 1397
$ws1.Range("F12").Value = 4644
$ws1.Range("F14").Value = 1750
$ws1.Range("F16").Value = 57
$ws1.Range("F18").Value = 204
$ws1.Range("F20").Value = 1045
$ws1.Range("F26").Value = 219
$ws1.Range("F30").Value = 105
$ws1.Range("F47").Value = 440

# Sheet "全部类型" (all types) - same events, slightly different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 1397
$ws4.Range("F12").Value = 4644
$ws4.Range("F13").Value = 1750
$ws4.Range("F16").Value = 57
$ws4.Range("F20").Value = 204
$ws4.Range("F23").Value = 1045
$ws4.Range("F26").Value = 219
$ws4.Range("F30").Value = 105
$ws4.Range("F44").Value = 440
